$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert new columns between Team (B) and ExpPoints ---
# New layout: A=Rank, B=Team, C=WIN, D=TOP4, E=TOP5, F=TOP6, G=RELEGATION, H=ExpPoints
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"

# copy the header style (bold, border, centered) from an existing header cell
# onto the new header cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("C1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Team re-ordering for ranks 11-13 (rows 12-14) ---
$ws.Range("B12").Value = "Valencia"
$ws.Range("B13").Value = "Osasuna"
$ws.Range("B14").Value = "Espanyol"

# --- New ExpPoints values (now living in column H) ---
$expPoints = @{
    2  = 87.79458757648719
    3  = 87.33773501937641
    4  = 68.6633083775025
    5  = 65.04670794206177
    6  = 61.61500156306601
    7  = 55.79356981534666
    8  = 54.25744093796727
    9  = 50.10166149086326
    10 = 49.8916889086379
    11 = 48.23327186422681
    12 = 47.09151727099622
    13 = 46.89074876220667
    14 = 46.40962468447881
    15 = 44.63937967515739
    16 = 40.58752684286147
    17 = 39.42187131853132
    18 = 38.82415045357858
    19 = 33.77498754230142
    20 = 33.33692311637653
    21 = 31.30720731208831
}

foreach ($row in 2..21) {
    # old ExpPoints value used to live in C; remove it and leave C:G blank
    # placeholders for the upcoming Monte Carlo simulation columns
    foreach ($col in @("C","D","E","F","G")) {
        $ws.Range("$col$row").Value = ""
        $ws.Range("$col$row").Style = "Normal"
    }
    $ws.Range("H$row").Value = $expPoints[$row]
}

$ws.Range("A1:H21").Columns.AutoFit() | Out-Null
